$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 21.84976866666667
$ws.Range("H2").Value = 65.549306
$ws.Range("I2").Value = 0.05020018890879543
$ws.Range("J2").Value = 0.05020018890879543
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 28.72417333333333
$ws.Range("N2").Value = 86.17251999999999
$ws.Range("O2").Value = 0.4233259107972328
$ws.Range("P2").Value = 0.4233259107972328
$ws.Range("Q2").Value = 627.6165424745687
$ws.Range("R2").Value = 5648.54888227112
$ws.Range("S2").Value = 0.02125104069200897
$ws.Range("T2").Value = 0.02125104069200897

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 21.84976866666667
$ws.Range("H3").Value = 65.549306
$ws.Range("I3").Value = 0.05020018890879543
$ws.Range("J3").Value = 0.05020018890879543
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 30.56986233333333
$ws.Range("N3").Value = 91.709587
$ws.Range("O3").Value = 0.4505269713084062
$ws.Range("P3").Value = 0.4505269713084062
$ws.Range("Q3").Value = 667.9444201551802
$ws.Range("R3").Value = 6011.499781396622
$ws.Range("S3").Value = 0.02261653906818945
$ws.Range("T3").Value = 0.02261653906818945

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 21.84976866666667
$ws.Range("H4").Value = 65.549306
$ws.Range("I4").Value = 0.05020018890879543
$ws.Range("J4").Value = 0.05020018890879543
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.559531999999999
$ws.Range("N4").Value = 25.678596
$ws.Range("O4").Value = 0.126147117894361
$ws.Range("P4").Value = 0.126147117894361
$ws.Range("Q4").Value = 187.0237940949307
$ws.Range("R4").Value = 1683.214146854376
$ws.Range("S4").Value = 0.006332609148597009
$ws.Range("T4").Value = 0.006332609148597009

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 385.0524703333334
$ws.Range("H5").Value = 1155.157411
$ws.Range("I5").Value = 0.8846641374295412
$ws.Range("J5").Value = 0.8846641374295412
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 28.72417333333333
$ws.Range("N5").Value = 86.17251999999999
$ws.Range("O5").Value = 0.4233259107972328
$ws.Range("P5").Value = 0.4233259107972328
$ws.Range("Q5").Value = 11060.31390028286
$ws.Range("R5").Value = 99542.82510254573
$ws.Range("S5").Value = 0.3745012517270088
$ws.Range("T5").Value = 0.3745012517270089

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 385.0524703333334
$ws.Range("H6").Value = 1155.157411
$ws.Range("I6").Value = 0.8846641374295412
$ws.Range("J6").Value = 0.8846641374295412
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 30.56986233333333
$ws.Range("N6").Value = 91.709587
$ws.Range("O6").Value = 0.4505269713084062
$ws.Range("P6").Value = 0.4505269713084062
$ws.Range("Q6").Value = 11771.00100919992
$ws.Range("R6").Value = 105939.0090827993
$ws.Range("S6").Value = 0.3985650544612949
$ws.Range("T6").Value = 0.3985650544612949

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 385.0524703333334
$ws.Range("H7").Value = 1155.157411
$ws.Range("I7").Value = 0.8846641374295412
$ws.Range("J7").Value = 0.8846641374295412
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.559531999999999
$ws.Range("N7").Value = 25.678596
$ws.Range("O7").Value = 0.126147117894361
$ws.Range("P7").Value = 0.126147117894361
$ws.Range("Q7").Value = 3295.868941497217
$ws.Range("R7").Value = 29662.82047347496
$ws.Range("S7").Value = 0.1115978312412375
$ws.Range("T7").Value = 0.1115978312412375

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 28.350479
$ws.Range("H8").Value = 85.05143699999999
$ws.Range("I8").Value = 0.06513567366166337
$ws.Range("J8").Value = 0.06513567366166337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 28.72417333333333
$ws.Range("N8").Value = 86.17251999999999
$ws.Range("O8").Value = 0.4233259107972328
$ws.Range("P8").Value = 0.4233259107972328
$ws.Range("Q8").Value = 814.3440728790265
$ws.Range("R8").Value = 7329.096655911238
$ws.Range("S8").Value = 0.02757361837821497
$ws.Range("T8").Value = 0.02757361837821497

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 28.350479
$ws.Range("H9").Value = 85.05143699999999
$ws.Range("I9").Value = 0.06513567366166337
$ws.Range("J9").Value = 0.06513567366166337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 30.56986233333333
$ws.Range("N9").Value = 91.709587
$ws.Range("O9").Value = 0.4505269713084062
$ws.Range("P9").Value = 0.4505269713084062
$ws.Range("Q9").Value = 866.6702401140575
$ws.Range("R9").Value = 7800.032161026518
$ws.Range("S9").Value = 0.02934537777892193
$ws.Range("T9").Value = 0.02934537777892193

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 28.350479
$ws.Range("H10").Value = 85.05143699999999
$ws.Range("I10").Value = 0.06513567366166337
$ws.Range("J10").Value = 0.06513567366166337
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.559531999999999
$ws.Range("N10").Value = 25.678596
$ws.Range("O10").Value = 0.126147117894361
$ws.Range("P10").Value = 0.126147117894361
$ws.Range("Q10").Value = 242.666832215828
$ws.Range("R10").Value = 2184.001489942452
$ws.Range("S10").Value = 0.008216677504526471
$ws.Range("T10").Value = 0.008216677504526471
